# oop/basic class diagrams: tweak text
#
# 1) Refresh the cached "Insert Date" auto-field (datetimeFigureOut) text on
#    every slide layout and on the slide master, 2/6/2017 -> 30/8/2017.
# 2) Nudge the little red textbox label on Slide 1 and retitle it from
#    "students" to "charges".
# 3) Re-create the two slide guides (best effort - older/limited hosts may
#    not persist these, so this step is not allowed to fail the script).

$p = $ppt.ActivePresentation

$oldDate = "2/6/2017"
$newDate = "30/8/2017"

# --- 1) Date placeholder on every custom layout -----------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $shapes = $layout.Shapes
    for ($si = 1; $si -le $shapes.Count; $si++) {
        $shp = $shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Date placeholder on the slide master itself -----------------------
$masterShapes = $master.Shapes
for ($si = 1; $si -le $masterShapes.Count; $si++) {
    $shp = $masterShapes.Item($si)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) "students" -> "charges" textbox on Slide 1 ----------------------
$slide = $p.Slides.Item(1)
$shapes1 = $slide.Shapes
for ($si = 1; $si -le $shapes1.Count; $si++) {
    $shp = $shapes1.Item($si)
    if ($shp.Name -eq "TextBox 93") {
        # EMU targets: off x 4191000 -> 4283968 (y stays 2971800),
        #              ext cx 1066800 -> 973832 (cy stays 369332).
        # Shape.Left/Width are points (Single-precision), so values are
        # chosen to land exactly on the target EMU after conversion.
        $shp.Left = 337.32035775878904
        $shp.Width = 76.67968503937009
        $shp.TextFrame.TextRange.Text = "charges"
    }
}

# --- 3) Slide guides (horizontal @2160, vertical @2880) -----------------
try {
    $guides = $p.Guides
    $g1 = $guides.Add(1, 2160)
    $g2 = $guides.Add(2, 2880)
} catch {
    # Guide manipulation isn't supported on every host; ignore failures.
}
